# Added information to the test plan:
#  - Developer name
#  - Method Inputs / Condition being Tested / Expected Result for the
#    8 ChequingAccount test cases (rows 7-14)
# The order in which the cells are written mirrors how the table was
# actually filled in (condition column first across several rows, then
# going back to fill the expected-result column), so that the resulting
# shared-string table is built up in the same sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$imports = "from bank_account.chequing_account import ChequingAccount`nfrom bank_account.bank_account import BankAccount`nimport unittest`nfrom datetime import date"

# Developer name
$ws.Range("C3").Value = "Ridham Sood"

# Row 7 - __init__ / appropriate value returned
$ws.Range("E7").Value = $imports
$ws.Range("F7").Value = "self.chequing = ChequingAccount(12345, 123, 1000, date(2025, 5, 1), 400, 0.02)"
$ws.Range("G7").Value = "Account Number: 12345`nClient Number: 123`nBalance: 1000`ndate_created: 2025-5-1`noverdraft_limit: 400`noverdraft_rate: 0.02"

# Rows 8-10 - conditions filled first
$ws.Range("E8").Value = $imports
$ws.Range("F8").Value = "account = ChequingAccount(12345, 123, 1000, date(2025, 5, 1), `"four hundred`", 0.02)`nexpected = -100"
$ws.Range("E9").Value = $imports
$ws.Range("F9").Value = "account = ChequingAccount(12345, 123, 1000, date(2025, 5, 1), 400, `"two`")`nexpected = 0.05"
$ws.Range("E10").Value = $imports
$ws.Range("F10").Value = "account = ChequingAccount(12345, 123, 1000, `"2025/5/1`", 400, 0.02)`nexpected = date.today()"

# ...then the expected results for rows 8-10
$ws.Range("G8").Value = "Overdraft Limit = -100"
$ws.Range("G9").Value = "Overdraft Rate = 0.05"
$ws.Range("G10").Value = "date_created = date.today()"

# Row 11
$ws.Range("E11").Value = $imports
$ws.Range("F11").Value = "account = ChequingAccount(12345, 123, 1000, date(2025, 5, 1), 500, 0.02)`nexpected = 0.50"

# Row 13 (filled before row 12's condition)
$ws.Range("E13").Value = $imports
$ws.Range("F13").Value = "account = ChequingAccount(12345, 123, 500, date(2025, 5, 1), 500, 0.02)`nexpected = 0.50"

# Row 12
$ws.Range("E12").Value = $imports
$ws.Range("F12").Value = "account = ChequingAccount(12345, 123, 200, date(2025, 5, 1), 500, 0.02)`nexpected = 6.50"

# Row 14
$ws.Range("E14").Value = $imports
$ws.Range("F14").Value = "account = ChequingAccount(12345, 123, 1000, date(2025, 5, 1), 500, 0.02)`nAccount number: 12345`nBalance: `$1,000.00`nOverdraft Limit: `$500.00`nOverdraft rate: `$2.00%`nAccount Type: Chequing"

# Expected results for rows 11, 12, 14, then 13 (reuses row 11's text)
$ws.Range("G11").Value = "Service Charge = 0.05"
$ws.Range("G12").Value = "Service Charge = 6.50"
$ws.Range("G14").Value = "Account number: 12345`nBalance: `$1,000.00`nOverdraft Limit: `$500.00`nOverdraft rate: `$2.00%`nAccount Type: Chequing"
$ws.Range("G13").Value = "Service Charge = 0.05"

# Leave the selection on G14, matching the saved view state
$ws.Range("G14").Select() | Out-Null
